$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 160, pushing existing rows 160-215 down to 161-216.
$ws.Rows.Item(160).Insert()

# Populate the newly inserted row 160 with the new daily record.
$ws.Cells.Item(160, 1).Value = 5
$ws.Cells.Item(160, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(160, 3).Value = "Maule"
$ws.Cells.Item(160, 4).Value = 44559
$ws.Cells.Item(160, 5).Value = 7
$ws.Cells.Item(160, 6).Value = 100112006
$ws.Cells.Item(160, 7).Value = "Repollo"
$ws.Cells.Item(160, 8).Value = "Crespo record"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 6000
$ws.Cells.Item(160, 11).Value = 600
$ws.Cells.Item(160, 12).Value = 600
$ws.Cells.Item(160, 13).Value = 600
$ws.Cells.Item(160, 14).Value = "$/unidad"
$ws.Cells.Item(160, 15).Value = "Región del Maule"
$ws.Cells.Item(160, 16).Value = 600
$ws.Cells.Item(160, 17).Value = 1
$ws.Cells.Item(160, 18).Value = "Hortaliza"
